# Daily attendance processing - 2026-02-07 01:54:36 UTC
# Update "Recorded By" column: swap "Administrator, Miss Dina Nasr"
# to "Miss Dina Nasr, Administrator" for every matching cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$oldText = "Administrator, Miss Dina Nasr"
$newText = "Miss Dina Nasr, Administrator"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
